$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "ec2-user" username with "ubuntu" for every credential row (7-32)
for ($r = 7; $r -le 32; $r++) {
    $ws.Range("F$r").Value = "ubuntu"
}

# Update the active selection to match the edited range
$ws.Range("F7:F32").Select()
